$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror the header cell's formatting (style) onto the new H1 header cell,
# same as the existing G1 "sum" header (bold, centered, thin border).
$ws.Range("G1").Copy($ws.Range("H1"))

# New "Save" column: header label + the row-2 flag value.
$ws.Range("H1").Value() = "Save"
$ws.Range("H2").Value() = 1
